# Generate Report for Handback
# The localization-status report is regenerated after a successful handback:
#  - the "Status" column flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" on every sheet that tracks it,
#  - the "Latest Handback DateTime" for each locale is refreshed to the
#    timestamp of the handback that just completed,
#  - the stale "Error Detail" (complaining the handback file was out of
#    date) is cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-07 05:03:57"
$zhcn.Range("P2").Value = ""

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-07 05:04:13"
$dede.Range("P2").Value = ""

# Resize the columns whose content width changed so the report reads
# cleanly (Status got longer, Error Detail got shorter/cleared).
$ovw.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ovw.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

$zhcn.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$zhcn.Columns.Item(16).EntireColumn.AutoFit() | Out-Null

$dede.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$dede.Columns.Item(16).EntireColumn.AutoFit() | Out-Null
